$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.042.35'
$ws.Range('E2').Value = '  -1.00%  '
$ws.Range('D3').Value = '1.826.71'
$ws.Range('E3').Value = '  -0.36%  '
$ws.Range('E4').Value = '  -0.30%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '310.57'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.38%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.007'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.32%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4622'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.46%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3660'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.77%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07248'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.60%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8611'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.85%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '19.90'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.94%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07814'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +6.62%  '
$ws.Range('D13').Value = '1.858.40'
$ws.Range('E13').Value = '  -1.04%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.336'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.76%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.545'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.25%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '91.80'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.007'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.30%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008684'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.26%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.008'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.20%  '
$ws.Range('B20').Value = 'WrappedBTC'
$ws.Range('C20').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D20').Value = '27.357.05'
$ws.Range('E20').Value = '  -0.63%  '
$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.51'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.84%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.160'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.36%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.55'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.03%  '
$ws.Range('D24').Value = '2.141.65'
$ws.Range('E24').Value = '  +2.33%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.22'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.11%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.843'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.63%  '
$ws.Range('E27').Value = '  -2.51%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.089'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.58%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.109'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.28%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '115.42'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.49%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08831'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.82%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.963'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.72%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.435'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.38%  '
$ws.Range('E34').Value = '  -3.71%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7208'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.082'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.16%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.443'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.51%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05239'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.01%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01936'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.89%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.948'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.90%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.207'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.33%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5162'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.56%  '
$ws.Range('B43').Value = 'Frax'
$ws.Range('C43').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8630'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -14.67%  '
$ws.Range('B44').Value = 'Algorand'
$ws.Range('C44').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1629'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.81%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.170'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.65%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4804'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.67%  '
$ws.Range('E47').Value = '  -0.37%  '
$ws.Range('E48').Value = '  -4.17%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '102.93'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.92%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06242'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.92%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.618'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.15%  '
